# Auto-generated edit script applying numeric updates to Chocobo_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76 (ALC)
$ws.Range("H76").Value = 3010.8147
$ws.Range("I76").Value = 2992.2917
$ws.Range("J76").Value = 3159
$ws.Range("K76").Value = 2992.2917
$ws.Range("L76").Value = 3159
$ws.Range("M76").Value = -2677.2917
$ws.Range("N76").Value = -3789

# Row 79 (ALC)
$ws.Range("H79").Value = 3010.8147
$ws.Range("I79").Value = 2992.2917
$ws.Range("J79").Value = 3159
$ws.Range("K79").Value = 2992.2917
$ws.Range("L79").Value = 3159
$ws.Range("M79").Value = -1900.2917
$ws.Range("N79").Value = -5343

# Row 112 (ALC)
$ws.Range("H112").Value = 1255.8195
$ws.Range("J112").Value = 1255.8195
$ws.Range("L112").Value = 3767.4585
$ws.Range("N112").Value = -5983.458500000001

# Row 137 (ALC)
$ws.Range("H137").Value = 901700.7
$ws.Range("I137").Value = 1834805.8
$ws.Range("K137").Value = 5504417.4
$ws.Range("M137").Value = -5501867.4

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws.Range("H61").Value = 1548.875
$ws.Range("I61").Value = 1513
$ws.Range("J61").Value = 1800
$ws.Range("K61").Value = 1513
$ws.Range("L61").Value = 1800
$ws.Range("M61").Value = -1301
$ws.Range("N61").Value = -2224

# Row 74 (ARM)
$ws.Range("H74").Value = 428039.12
$ws.Range("I74").Value = 592910.75
$ws.Range("J74").Value = 4083.4285
$ws.Range("K74").Value = 592910.75
$ws.Range("L74").Value = 4083.4285
$ws.Range("M74").Value = -592036.75
$ws.Range("N74").Value = -5831.4285

# Row 77 (ARM)
$ws.Range("H77").Value = 428039.12
$ws.Range("I77").Value = 592910.75
$ws.Range("J77").Value = 4083.4285
$ws.Range("K77").Value = 2964553.75
$ws.Range("L77").Value = 20417.1425
$ws.Range("M77").Value = -2960185.75
$ws.Range("N77").Value = -29153.1425

# Row 97 (ARM)
$ws.Range("H97").Value = 648.0454999999999
$ws.Range("I97").Value = 516.0625
$ws.Range("K97").Value = 516.0625
$ws.Range("M97").Value = -20.0625

# Row 110 (ARM)
$ws.Range("H110").Value = 2196.25
$ws.Range("I110").Value = 1756
$ws.Range("J110").Value = 2930
$ws.Range("K110").Value = 1756
$ws.Range("L110").Value = 2930
$ws.Range("M110").Value = 289
$ws.Range("N110").Value = -7020

# Row 122 (ARM)
$ws.Range("H122").Value = 3259
$ws.Range("I122").Value = 1448
$ws.Range("J122").Value = 5070
$ws.Range("K122").Value = 4344
$ws.Range("L122").Value = 15210
$ws.Range("M122").Value = -1894
$ws.Range("N122").Value = -20110

# Row 132 (ARM)
$ws.Range("H132").Value = 2662.375
$ws.Range("I132").Value = 1471.3334
$ws.Range("J132").Value = 4448.9375
$ws.Range("K132").Value = 4414.0002
$ws.Range("L132").Value = 13346.8125
$ws.Range("M132").Value = -1884.0002
$ws.Range("N132").Value = -18406.8125

# Row 136 (ARM)
$ws.Range("H136").Value = 1548.875
$ws.Range("I136").Value = 1513
$ws.Range("J136").Value = 1800
$ws.Range("K136").Value = 4539
$ws.Range("L136").Value = 5400
$ws.Range("M136").Value = -1989
$ws.Range("N136").Value = -10500

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (BSM)
$ws.Range("H20").Value = 8103.174
$ws.Range("I20").Value = 1236.909
$ws.Range("J20").Value = 14397.25
$ws.Range("K20").Value = 1236.909
$ws.Range("L20").Value = 14397.25
$ws.Range("M20").Value = -989.9090000000001
$ws.Range("N20").Value = -14891.25

# Row 22 (BSM)
$ws.Range("H22").Value = 186.75
$ws.Range("I22").Value = 186.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 186.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -13.75
$ws.Range("N22").ClearContents()

# Row 62 (BSM)
$ws.Range("H62").Value = 29000
$ws.Range("J62").Value = 29000
$ws.Range("L62").Value = 29000
$ws.Range("N62").Value = -30372

# Row 65 (BSM)
$ws.Range("H65").Value = 29000
$ws.Range("J65").Value = 29000
$ws.Range("L65").Value = 87000
$ws.Range("N65").Value = -93864

# Row 99 (BSM)
$ws.Range("H99").Value = 9853.154
$ws.Range("I99").Value = 2172.5
$ws.Range("K99").Value = 2172.5
$ws.Range("M99").Value = -674.5

# Row 134 (BSM)
$ws.Range("H134").Value = 4261.0557
$ws.Range("I134").Value = 2053.1428
$ws.Range("J134").Value = 5666.091
$ws.Range("K134").Value = 6159.428400000001
$ws.Range("L134").Value = 16998.273
$ws.Range("M134").Value = -3624.428400000001
$ws.Range("N134").Value = -22068.273

# Row 140 (BSM)
$ws.Range("H140").Value = 45289.633
$ws.Range("J140").Value = 45289.633
$ws.Range("L140").Value = 45289.633
$ws.Range("N140").Value = -55649.633

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 267448.2
$ws.Range("I31").Value = 712318.9399999999
$ws.Range("J31").Value = 3306.1562
$ws.Range("K31").Value = 712318.9399999999
$ws.Range("L31").Value = 3306.1562
$ws.Range("M31").Value = -712023.9399999999
$ws.Range("N31").Value = -3896.1562

# Row 34 (CRP)
$ws.Range("H34").Value = 267448.2
$ws.Range("I34").Value = 712318.9399999999
$ws.Range("J34").Value = 3306.1562
$ws.Range("K34").Value = 712318.9399999999
$ws.Range("L34").Value = 3306.1562
$ws.Range("M34").Value = -712116.9399999999
$ws.Range("N34").Value = -3710.1562

# Row 52 (CRP)
$ws.Range("H52").Value = 43814.285
$ws.Range("I52").Value = 15000
$ws.Range("J52").Value = 55340
$ws.Range("K52").Value = 15000
$ws.Range("L52").Value = 55340
$ws.Range("M52").Value = -14706
$ws.Range("N52").Value = -55928

# Row 99 (CRP)
$ws.Range("H99").Value = 11768507
$ws.Range("I99").Value = 20001570
$ws.Range("J99").Value = 6987.7144
$ws.Range("K99").Value = 20001570
$ws.Range("L99").Value = 6987.7144
$ws.Range("M99").Value = -20000072
$ws.Range("N99").Value = -9983.714400000001

# Row 126 (CRP)
$ws.Range("H126").Value = 11768507
$ws.Range("I126").Value = 20001570
$ws.Range("J126").Value = 6987.7144
$ws.Range("K126").Value = 60004710
$ws.Range("L126").Value = 20963.1432
$ws.Range("M126").Value = -60002240
$ws.Range("N126").Value = -25903.1432

# Row 132 (CRP)
$ws.Range("H132").Value = 5308.1333
$ws.Range("I132").Value = 4762.4
$ws.Range("J132").Value = 6399.6
$ws.Range("K132").Value = 14287.2
$ws.Range("L132").Value = 19198.8
$ws.Range("M132").Value = -11757.2
$ws.Range("N132").Value = -24258.8

# Row 134 (CRP)
$ws.Range("H134").Value = 2482.55
$ws.Range("I134").Value = 1038.4117
$ws.Range("J134").Value = 10666
$ws.Range("K134").Value = 3115.2351
$ws.Range("L134").Value = 31998
$ws.Range("M134").Value = -580.2351000000003
$ws.Range("N134").Value = -37068

$ws = $wb.Worksheets.Item("CUL")
# Row 113 (CUL)
$ws.Range("H113").Value = 1061.5454
$ws.Range("J113").Value = 2226.3333
$ws.Range("L113").Value = 6678.999899999999
$ws.Range("N113").Value = -11018.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 63 (GSM)
$ws.Range("H63").Value = 19566.5
$ws.Range("J63").Value = 19566.5
$ws.Range("L63").Value = 19566.5
$ws.Range("N63").Value = -20938.5

# Row 66 (GSM)
$ws.Range("H66").Value = 19566.5
$ws.Range("J66").Value = 19566.5
$ws.Range("L66").Value = 58699.5
$ws.Range("N66").Value = -65563.5

# Row 70 (GSM)
$ws.Range("H70").Value = 6716.6
$ws.Range("I70").Value = 6369
$ws.Range("J70").Value = 7411.8
$ws.Range("K70").Value = 6369
$ws.Range("L70").Value = 7411.8
$ws.Range("M70").Value = -6099
$ws.Range("N70").Value = -7951.8

# Row 73 (GSM)
$ws.Range("H73").Value = 6716.6
$ws.Range("I73").Value = 6369
$ws.Range("J73").Value = 7411.8
$ws.Range("K73").Value = 6369
$ws.Range("L73").Value = 7411.8
$ws.Range("M73").Value = -5433
$ws.Range("N73").Value = -9283.799999999999

# Row 113 (GSM)
$ws.Range("H113").Value = 1224
$ws.Range("I113").Value = 1130.5
$ws.Range("J113").Value = 1348.6666
$ws.Range("K113").Value = 1130.5
$ws.Range("L113").Value = 1348.6666
$ws.Range("M113").Value = 1039.5
$ws.Range("N113").Value = -5688.6666

# Row 126 (GSM)
$ws.Range("H126").Value = 3734.415
$ws.Range("I126").Value = 2754.027
$ws.Range("K126").Value = 8262.081
$ws.Range("M126").Value = -5792.081

# Row 132 (GSM)
$ws.Range("H132").Value = 3437.3416
$ws.Range("I132").Value = 2339.8635
$ws.Range("J132").Value = 4708.1055
$ws.Range("K132").Value = 7019.5905
$ws.Range("L132").Value = 14124.3165
$ws.Range("M132").Value = -4489.5905
$ws.Range("N132").Value = -19184.3165

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (LTW)
$ws.Range("H40").Value = 5106.7085
$ws.Range("I40").Value = 4538.5884
$ws.Range("J40").Value = 6486.4287
$ws.Range("K40").Value = 4538.5884
$ws.Range("L40").Value = 6486.4287
$ws.Range("M40").Value = -4402.5884
$ws.Range("N40").Value = -6758.4287

# Row 132 (LTW)
$ws.Range("H132").Value = 6860.4
$ws.Range("I132").Value = 5240.6
$ws.Range("J132").Value = 7670.3
$ws.Range("K132").Value = 15721.8
$ws.Range("L132").Value = 23010.9
$ws.Range("M132").Value = -13191.8
$ws.Range("N132").Value = -28070.9

# Row 136 (LTW)
$ws.Range("H136").Value = 3290.4119
$ws.Range("I136").Value = 937.9048
$ws.Range("J136").Value = 7090.615
$ws.Range("K136").Value = 2813.7144
$ws.Range("L136").Value = 21271.845
$ws.Range("M136").Value = -263.7143999999998
$ws.Range("N136").Value = -26371.845

$ws = $wb.Worksheets.Item("WVR")
# Row 46 (WVR)
$ws.Range("H46").Value = 70074
$ws.Range("J46").Value = 70074
$ws.Range("L46").Value = 70074
$ws.Range("N46").Value = -70536

# Row 108 (WVR)
$ws.Range("H108").Value = 29450
$ws.Range("J108").Value = 29450
$ws.Range("L108").Value = 29450
$ws.Range("N108").Value = -37130

# Row 132 (WVR)
$ws.Range("H132").Value = 7411917
$ws.Range("I132").Value = 7421.5
$ws.Range("J132").Value = 11497156
$ws.Range("K132").Value = 22264.5
$ws.Range("L132").Value = 34491468
$ws.Range("M132").Value = -19734.5
$ws.Range("N132").Value = -34496528

# Row 134 (WVR)
$ws.Range("H134").Value = 70074
$ws.Range("J134").Value = 70074
$ws.Range("L134").Value = 210222
$ws.Range("N134").Value = -215292
